# Apply the edits described by the diff to the PCQA/ExcelDatosCuentas.xlsx workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update data rows 2 and 3 (B, C, F, G, H, O columns) ---
# Row 2 (B2 keeps its quote-prefixed "text" formatting, hence the leading apostrophe)
$ws.Range("B2").Value = "'preproducciongestion.segurossura.com.ar"
$ws.Range("C2").Value = "https://preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("G2").Value = "PruebaRegre"
$ws.Range("F2").Value = "AnswerRegre"
$ws.Range("H2").Value = 20300114
$ws.Range("O2").Value = 118

# Row 3
$ws.Range("B3").Value = "'preproducciongestion.segurossura.com.ar"
$ws.Range("C3").Value = "https://preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("F3").Value = "MattioliRegre"
$ws.Range("G3").Value = "PruebaRegre"
$ws.Range("H3").Value = 20300115
$ws.Range("O3").Value = 119

# --- Update sheet view (scroll position / selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("O4").Select()
